$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 3) with a second quote for testing random pull
# Set C3's quote text first so the shared-string table order matches
# (quote string gets added before author string).
$ws.Range("A3").Value = "habits"
$ws.Range("C3").Value = "Repetition of the same thought or physical action develops into a habit which, repeated frequently enough, becomes an automatic reflex."
$ws.Range("B3").Value = "Norman Vincent Peale"

# Match formatting of row 2 (wrap text on column C, same row height)
$ws.Range("C3").WrapText = $true
$ws.Rows.Item(3).RowHeight = $ws.Rows.Item(2).RowHeight

# Update selection as in diff
$ws.Range("B8").Select()
